$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several "Price" (column D) values look like plain decimals (e.g. "43.12").
# Setting .Value directly on those would let Excel auto-convert them to numbers
# (losing the text formatting / introducing float rounding, e.g. 43.119999999999997).
# To keep them as literal text - matching the original inline-string cells - we
# write a text-producing formula, then Copy + PasteSpecial(xlPasteValues, -4163)
# to flatten it back down to a plain text value without allocating a new
# NumberFormat style (avoiding any styles.xml changes).

$ws.Range("D2").Value = "69.519.79"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "3.943.11"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Formula = '="492.25"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Formula = '="146.84"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +5.42%  "
$ws.Range("D11").Formula = '="0.0000350"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("D12").Formula = '="43.12"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Formula = '="10.48"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "4.572.12"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "3.958.39"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Formula = '="14.27"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -4.21%  "
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").Value = "69.507.07"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").Formula = '="441.17"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  +1.70%  "
$ws.Range("D23").Formula = '="14.53"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").Formula = '="89.37"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").Formula = '="12.00"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +9.46%  "
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").Formula = '="11.12"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -4.38%  "
$ws.Range("D28").Formula = '="37.30"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -4.66%  "
$ws.Range("D29").Formula = '="5.66"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -3.58%  "
$ws.Range("D30").Formula = '="707.18"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Formula = '="13.52"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Formula = '="0.475"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +26.74%  "
$ws.Range("D35").Value = "0.0₃0908"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("D36").Formula = '="61.67"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D37").Formula = '="6.07"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +4.43%  "
$ws.Range("D38").Formula = '="40.82"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +2.55%  "
$ws.Range("D43").Formula = '="2.93"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Formula = '="0.144"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("E47").Value = "  +10.23%  "
$ws.Range("D48").Formula = '="3.35"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +6.50%  "
$ws.Range("E49").Value = "  +8.26%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").Formula = '="144.04"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -1.51%  "

$excel.CutCopyMode = 0

